$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.172.76"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.601.82"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "303.39"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("D7").Value = "0.3782"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "52.08"
$ws.Range("E8").Value = "  +4.38%  "
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").Value = "1.270"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "0.08124"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "22.80"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").Value = "6.596"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "7.419"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "0.00001245"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").Value = "1.602.48"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "94.04"
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("D19").Value = "0.06880"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").Value = "6.548"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").Value = "23.174.82"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "2.399"
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("D26").Value = "2.973"
$ws.Range("E26").Value = "  +7.55%  "
$ws.Range("D27").Value = "21.26"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "149.27"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").Value = "133.99"
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").Value = "2.361"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "6.776"
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("D33").Value = "1.781.30"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").Value = "0.9702"
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("D35").Value = "0.07523"
$ws.Range("E35").Value = "  -2.44%  "
$ws.Range("D36").Value = "10.27"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("D37").Value = "0.02720"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "0.2504"
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D40").Value = "6.073"
$ws.Range("E40").Value = "  -3.39%  "
$ws.Range("D41").Value = "0.7108"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "1.362"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").Value = "12.50"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("D44").Value = "15.68"
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("D45").Value = "0.6535"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").Value = "4.014"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "132.27"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").Value = "0.07950"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("E51").Value = "  +1.04%  "
